$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fix the Treatment Tab query text in B5: drop the redundant CONCAT()
# wrapper around REPLACE(trt.treatment_agent, ';', ', ') ---------------------
$oldAgentExpr = "CONCAT(REPLACE(trt.treatment_agent, ';', ', ')) AS ""Treatment Agent"","
$newAgentExpr = "REPLACE(trt.treatment_agent, ';', ', ') AS ""Treatment Agent"","
$treatmentQuery = $ws.Range("B5").Value2
$treatmentQuery = $treatmentQuery.Replace($oldAgentExpr, $newAgentExpr)
$ws.Range("B5").Value = $treatmentQuery

# --- Match the font/wrap formatting already used by the Treatment Tab query
# cell (B5) on the other TabQuery cells that were restyled (B2, B3, B4) -----
foreach ($addr in @("B2", "B3", "B4")) {
    $cell = $ws.Range($addr)
    $cell.Font.Size = 12
    $cell.WrapText = $true
}

# --- Move the visible selection to C5 (scrolled down to show row 4 onward) -
$ws.Range("C5").Select()
